$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("covariate_importance")
$ws1.Range("A2").Value = "state"
$ws1.Range("B2").Value = 100
$ws1.Range("C2").Value = 100
$ws1.Range("D2").Value = 100
$ws1.Range("A3").Value = "region"
$ws1.Range("B3").Value = 100
$ws1.Range("C3").Value = 98
$ws1.Range("D3").Value = 98
$ws1.Range("A4").Value = "percenttwoormoreraces"
$ws1.Range("B4").Value = 73
$ws1.Range("C4").Value = 46
$ws1.Range("D4").Value = 46
$ws1.Range("A5").Value = "percentfreelunchqualified"
$ws1.Range("B5").Value = 61
$ws1.Range("C5").Value = 40
$ws1.Range("D5").Value = 40
$ws1.Range("A6").Value = "cntycaseschange"
$ws1.Range("B6").Value = 65
$ws1.Range("C6").Value = 34
$ws1.Range("D6").Value = 34
$ws1.Range("A7").Value = "percentwhite"
$ws1.Range("B7").Value = 51
$ws1.Range("C7").Value = 30
$ws1.Range("D7").Value = 30
$ws1.Range("A8").Value = "percentamericanindianoralaskanative"
$ws1.Range("B8").Value = 55.00000000000001
$ws1.Range("C8").Value = 27
$ws1.Range("D8").Value = 27
$ws1.Range("A9").Value = "rplthemes"
$ws1.Range("B9").Value = 46
$ws1.Range("C9").Value = 26
$ws1.Range("D9").Value = 26
$ws1.Range("A10").Value = "percentblackorafricanamerican"
$ws1.Range("B10").Value = 46
$ws1.Range("C10").Value = 22
$ws1.Range("D10").Value = 22
$ws1.Range("A11").Value = "percentnativehawaiianorotherpacificislander"
$ws1.Range("B11").Value = 56.99999999999999
$ws1.Range("C11").Value = 17
$ws1.Range("D11").Value = 17
$ws1.Range("A12").Value = "percentasian"
$ws1.Range("B12").Value = 34
$ws1.Range("C12").Value = 15
$ws1.Range("D12").Value = 15
$ws1.Range("A13").Value = "derivedtotalenrolled"
$ws1.Range("B13").Value = 51
$ws1.Range("C13").Value = 13
$ws1.Range("D13").Value = 13
$ws1.Range("A14").Value = "locale"
$ws1.Range("B14").Value = 35
$ws1.Range("C14").Value = 13
$ws1.Range("D14").Value = 13
$ws1.Range("A15").Value = "schoollevel"
$ws1.Range("B15").Value = 30
$ws1.Range("C15").Value = 9
$ws1.Range("D15").Value = 9
$ws1.Range("A16").Value = "percenthispaniclatino"
$ws1.Range("B16").Value = 34
$ws1.Range("C16").Value = 7.000000000000001
$ws1.Range("D16").Value = 7.000000000000001
$ws1.Range("A17").Value = "percentnotspecified"
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 3
$ws1.Range("D17").Value = 0

$ws2 = $wb.Worksheets.Item("strategy_importance")
$ws2.Range("A2").Value = "contacttracing"
$ws2.Range("B2").Value = 98
$ws2.Range("C2").Value = 98
$ws2.Range("D2").Value = 98
$ws2.Range("A3").Value = "hvacsystems"
$ws2.Range("B3").Value = 98
$ws2.Range("C3").Value = 99
$ws2.Range("D3").Value = 98
$ws2.Range("A4").Value = "cleaning"
$ws2.Range("B4").Value = 80
$ws2.Range("C4").Value = 78
$ws2.Range("D4").Value = 77
$ws2.Range("A5").Value = "screeningtestingforstudents"
$ws2.Range("B5").Value = 59
$ws2.Range("C5").Value = 51
$ws2.Range("D5").Value = 56.99999999999999
$ws2.Range("A6").Value = "quarantine"
$ws2.Range("B6").Value = 26
$ws2.Range("C6").Value = 18
$ws2.Range("D6").Value = 25
$ws2.Range("A7").Value = "masks"
$ws2.Range("B7").Value = 20
$ws2.Range("C7").Value = 15
$ws2.Range("D7").Value = 16
$ws2.Range("A8").Value = "physicaldistancing"
$ws2.Range("B8").Value = 15
$ws2.Range("C8").Value = 10
$ws2.Range("D8").Value = 14
$ws2.Range("A9").Value = "vaccination"
$ws2.Range("B9").Value = 17
$ws2.Range("C9").Value = 17
$ws2.Range("D9").Value = 14
$ws2.Range("A10").Value = "hepafilters"
$ws2.Range("B10").Value = 14
$ws2.Range("C10").Value = 12
$ws2.Range("D10").Value = 8
$ws2.Range("A11").Value = "stayhome"
$ws2.Range("B11").Value = 4
$ws2.Range("C11").Value = 2
$ws2.Range("D11").Value = 3
